$wb = $excel.ActiveWorkbook

$oldDataSheetName = "Durable acc. 1995-2019"
$newDataSheetName = "Figure 3"

# --- Remove the now-unnecessary sheet and keep only the data/figure sheet ---
# "Durable Accbysector&prod 2011" held a different, now-unused breakdown;
# only "Durable acc. 1995-2019" is needed to reproduce Figure 3.
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Durable Accbysector&prod 2011").Delete()
$excel.DisplayAlerts = $true

# --- Rename the remaining sheet ---
$ws = $wb.Worksheets.Item($oldDataSheetName)
$ws.Name = $newDataSheetName

# --- Fix up the chart's series formulas so they point at the renamed sheet ---
# (the engine doesn't auto-propagate a sheet rename into chart series refs)
foreach ($co in $ws.ChartObjects()) {
    $chart = $co.Chart
    $series = $chart.SeriesCollection()
    for ($i = 1; $i -le $series.Count; $i++) {
        $s = $series.Item($i)
        $s.Formula = $s.Formula -replace [regex]::Escape("'" + $oldDataSheetName + "'"), ("'" + $newDataSheetName + "'")
    }
}

$ws.Activate()
